# Week 3 day 5
# Update the title slide's headline from "Week 2 Recap" to "Week 3 Recap".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Week 3 Recap"
